# Rename "exac" to "pop_stats" everywhere it appears as a cell value.
# (Test Cases List sheet, column D, rows 47-56 currently hold the literal
# string "exac" -- the commit renames that value to "pop_stats".)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -eq "exac") {
            $cell.Value = "pop_stats"
        }
    }
}

# Move the active selection, matching the recorded cursor position after the edit.
$null = $ws.Range("M43").Select()
